$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 66, shifting existing rows 66-110 down to 67-111
$ws.Rows.Item(66).Insert()

# Populate the newly inserted row 66 with the new record
$ws.Cells.Item(66, 1).Value = 5
$ws.Cells.Item(66, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(66, 3).Value = "Maule"
$ws.Cells.Item(66, 4).Value = 44606
$ws.Cells.Item(66, 5).Value = 7
$ws.Cells.Item(66, 6).Value = 100112030
$ws.Cells.Item(66, 7).Value = "Poroto granado"
$ws.Cells.Item(66, 8).Value = "Sin especificar"
$ws.Cells.Item(66, 9).Value = "Primera"
$ws.Cells.Item(66, 10).Value = 300
$ws.Cells.Item(66, 11).Value = 20000
$ws.Cells.Item(66, 12).Value = 20000
$ws.Cells.Item(66, 13).Value = 20000
$ws.Cells.Item(66, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(66, 15).Value = "Región del Maule"
$ws.Cells.Item(66, 16).Value = 800
$ws.Cells.Item(66, 17).Value = 25
$ws.Cells.Item(66, 18).Value = "Hortaliza"

# Apply the same date cell number format used by the other date cells in column D
$ws.Cells.Item(66, 4).NumberFormat = $ws.Cells.Item(67, 4).NumberFormat
